# Add a new appendix row ("Using Digital Multimeters") to the "lablist"
# worksheet, inserted as a new row 73 (pushing the former rows 73-88, and
# the summary rows below, down by one row), per the commit:
#   "New Appendix on using Digital Multimeters"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lablist")

# --- 1. Insert a new row at 73; Excel shifts all formulas/rows below it
#        down automatically (incl. the SUMIF($V$2:V$79...) ranges etc.) ---
$ws.Rows("73:73").Insert()

# --- 2. Copy formatting for the new row 73 from matching donor cells
#        elsewhere in the sheet that already carry the desired style ---
$ws.Range("A8").Copy()
$ws.Range("A73").PasteSpecial(-4122)

$ws.Range("B72").Copy()
$ws.Range("B73").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("C73").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("D73:J73").PasteSpecial(-4122)

$ws.Range("K5").Copy()
$ws.Range("K73").PasteSpecial(-4122)

$ws.Range("P4").Copy()
$ws.Range("L73:U73").PasteSpecial(-4122)

$ws.Range("V4").Copy()
$ws.Range("V73").PasteSpecial(-4122)

$ws.Range("X35").Copy()
$ws.Range("X73").PasteSpecial(-4122)

# --- 3. Fill in the new row's content ---
$ws.Range("B73").Value = "Using Digital Multimeters"
$ws.Range("K73").Value = "New in fall 2019, by MT"
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 1
$ws.Range("F73").Formula = "=CEILING(D73,X`$22+1)"

# --- 4. The row-insert operation loses the CSE "array formula" marker on
#        the summary block below (rows 81/82/83/86 after the shift); the
#        formula text is fine but it must be re-entered as an array
#        formula so it evaluates (SUM of an array product) instead of
#        erroring out with #VALUE!. Re-apply FormulaArray on each. ---
$cols = @("L","M","N","O","P","Q")
foreach ($c in $cols) {
  $ws.Range($c + "81").FormulaArray = "=SUM(`$F2:`$F79*(" + $c + "2:" + $c + "79>=0.9)*(`$V2:`$V79>=`$X`$12))"
  $ws.Range($c + "82").FormulaArray = "=SUM(`$F2:`$F79*" + $c + "2:" + $c + "79*(`$V2:`$V79>=`$X`$12))"
  $ws.Range($c + "83").FormulaArray = "=SUM(`$F`$2:`$F`$79*(" + $c + "`$2:" + $c + "`$79>=0.1)*(`$V`$2:`$V`$79>=`$X`$12))"
  $ws.Range($c + "86").FormulaArray = "=SUM(`$F2:`$F79*" + $c + "2:" + $c + "79*(`$V2:`$V79<`$X`$12))"
}

# --- 5. Recalculate everything so cached values are consistent. ---
$excel.Calculate()
